$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Header date: "March 7, 2025" -> "{dateRequested}"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("March 7, 2025", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{dateRequested}", 1) | Out-Null

# ---------------------------------------------------------------------------
# 2) "___________________, is presently residing at ____________________
#     personally appeared in this office    and declared that he/she has no
#     sufficient source of income to provide funds for his/her studies"
#    -> "{fullName}, is presently residing at {address} personally appeared
#        in this office and declared that he/she has  no sufficient source
#        of income to provide funds for his/her studies"
# ---------------------------------------------------------------------------

# 2a) the 19-underscore blank -> "{fullName},"  (plain for now, underline next)
$d.Content.Find.Execute("___________________", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{fullName},", 1) | Out-Null

# 2b) underline just the word "fullName" (splits the merged run automatically)
$rngFullName = $d.Content
$rngFullName.Find.Execute("fullName") | Out-Null
$rngFullName.Font.Underline = 1

# 2c) drop the now-duplicated leading comma on the next run
$d.Content.Find.Execute(", is presently residing at ", $false, $false, $false, $false, $false, `
    $true, 1, $false, " is presently residing at ", 1) | Out-Null

# 2d) the 20-underscore blank (+ trailing space) -> "{address} "
$d.Content.Find.Execute("____________________ ", $false, $false, $false, $false, $false, `
    $true, 1, $false, "{address} ", 1) | Out-Null

# 2e) underline just the word "address"
$rngAddress = $d.Content
$rngAddress.Find.Execute("address") | Out-Null
$rngAddress.Font.Underline = 1

# 2f) tidy up the spacing in the trailing sentence
$d.Content.Find.Execute("personally appeared in this office    and declared that he/she has no sufficient source of income to provide funds for his/her studies", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "personally appeared in this office and declared that he/she has  no sufficient source of income to provide funds for his/her studies", 1) | Out-Null

# ---------------------------------------------------------------------------
# 3) " is issued upon the request for ______________________ purposes"
#    -> " is issued upon the request for    {purpose} purposes"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("______________________ purposes", $false, $false, $false, $false, $false, `
    $true, 1, $false, "   {purpose} purposes", 1) | Out-Null

$rngPurpose = $d.Content
$rngPurpose.Find.Execute("purpose") | Out-Null
$rngPurpose.Font.Underline = 1

# ---------------------------------------------------------------------------
# 4) "This certificate is valid until October 2025." -> "... {validUntil}."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This certificate is valid until October 2025.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "This certificate is valid until {validUntil}.", 1) | Out-Null

Write-Host "edits applied"
